$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.01725977985153514
$ws.Range("C2").Value = 0.2126533318958536

$ws.Range("B3").Value = 0.06494109381072437
$ws.Range("C3").Value = 0.2253125392105234

$ws.Range("B4").Value = 0.8704471580490418
$ws.Range("C4").Value = 0.5474270679820905

$ws.Range("B5").Value = 0.9933034061575118
$ws.Range("C5").Value = 0.4673423134879546

$ws.Range("B6").Value = 0.9716928825097153
$ws.Range("C6").Value = 0.8067378392971114

$ws.Range("B7").Value = 0.9459351052397782
$ws.Range("C7").Value = 0.3767165195422588

$ws.Range("B8").Value = 0.008739812970161439
$ws.Range("C8").Value = 0.1716621780395508

$ws.Range("B9").Value = 0.1226151811912138
$ws.Range("C9").Value = 0.2223454458028964

$ws.Range("B10").Value = 0.6911459381654664
$ws.Range("C10").Value = 0.3483966465056509
